$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = "Colours"
$ws.Range("B8").Value = "Basic asoociation recogniceble patterns"

$ws.Range("A9").Select()
